# Auto-generated edit script: updates cached market/profit values
# on the ALC, BSM, CRP, CUL, GSM, LTW, WVR sheets to match the
# scheduled runner's latest data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H62").Value = 8614.214
$ws.Range("I62").Value = 4405.5
$ws.Range("J62").Value = 10297.7
$ws.Range("K62").Value = 4405.5
$ws.Range("L62").Value = 10297.7
$ws.Range("M62").Value = -3781.5
$ws.Range("N62").Value = -11545.7
$ws.Range("H64").Value = 7803.3335
$ws.Range("I64").Value = 6490.2
$ws.Range("J64").Value = 8459.9
$ws.Range("K64").Value = 6490.2
$ws.Range("L64").Value = 8459.9
$ws.Range("M64").Value = -6242.2
$ws.Range("N64").Value = -8955.9
$ws.Range("H65").Value = 8614.214
$ws.Range("I65").Value = 4405.5
$ws.Range("J65").Value = 10297.7
$ws.Range("K65").Value = 22027.5
$ws.Range("L65").Value = 51488.5
$ws.Range("M65").Value = -18907.5
$ws.Range("N65").Value = -57728.5
$ws.Range("H67").Value = 7803.3335
$ws.Range("I67").Value = 6490.2
$ws.Range("J67").Value = 8459.9
$ws.Range("K67").Value = 6490.2
$ws.Range("L67").Value = 8459.9
$ws.Range("M67").Value = -5632.2
$ws.Range("N67").Value = -10175.9
$ws.Range("H86").Value = 2541.8333
$ws.Range("I86").Value = 2250.25
$ws.Range("K86").Value = 2250.25
$ws.Range("M86").Value = -1127.25
$ws.Range("H89").Value = 2541.8333
$ws.Range("I89").Value = 2250.25
$ws.Range("K89").Value = 11251.25
$ws.Range("M89").Value = -5635.25
$ws.Range("H112").Value = 1091.75
$ws.Range("I112").Value = 1213.6666
$ws.Range("J112").Value = 726
$ws.Range("K112").Value = 3640.9998
$ws.Range("L112").Value = 2178
$ws.Range("M112").Value = -2532.9998
$ws.Range("N112").Value = -4394
$ws.Range("H115").Value = 740.75
$ws.Range("I115").Value = 740.75
$ws.Range("K115").Value = 2222.25
$ws.Range("M115").Value = -655.25
$ws.Range("H131").Value = 1558.8334
$ws.Range("I131").Value = 1470.6
$ws.Range("K131").Value = 4411.799999999999
$ws.Range("M131").Value = 628.2000000000007
$ws.Range("H132").Value = 22096.834
$ws.Range("I132").Value = 29823.143
$ws.Range("K132").Value = 89469.429
$ws.Range("M132").Value = -86939.429
$ws.Range("H137").Value = 2404.3076
$ws.Range("I137").Value = 1440.6666
$ws.Range("J137").Value = 4572.5
$ws.Range("K137").Value = 4321.9998
$ws.Range("L137").Value = 13717.5
$ws.Range("M137").Value = -1771.9998
$ws.Range("N137").Value = -18817.5

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H107").Value = 4185.5835
$ws.Range("I107").Value = 1342.6923
$ws.Range("K107").Value = 1342.6923
$ws.Range("M107").Value = 577.3077000000001

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H86").Value = 2499.5
$ws.Range("I86").Value = 2499.5
$ws.Range("K86").Value = 2499.5
$ws.Range("M86").Value = -1376.5
$ws.Range("H89").Value = 2499.5
$ws.Range("I89").Value = 2499.5
$ws.Range("K89").Value = 12497.5
$ws.Range("M89").Value = -6881.5
$ws.Range("H97").Value = 30000
$ws.Range("J97").Value = 30000
$ws.Range("L97").Value = 30000
$ws.Range("N97").Value = -31982
$ws.Range("H99").Value = 4087.3
$ws.Range("I99").Value = 3814
$ws.Range("K99").Value = 3814
$ws.Range("M99").Value = -2316
$ws.Range("H126").Value = 4087.3
$ws.Range("I126").Value = 3814
$ws.Range("K126").Value = 11442
$ws.Range("M126").Value = -8972
$ws.Range("H132").Value = 4561.2573
$ws.Range("I132").Value = 3671.762
$ws.Range("K132").Value = 11015.286
$ws.Range("M132").Value = -8485.286
$ws.Range("H134").Value = 1442.9333
$ws.Range("I134").Value = 1442.9333
$ws.Range("K134").Value = 4328.7999
$ws.Range("M134").Value = -1793.7999

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H132").Value = 1899.25
$ws.Range("I132").Value = 800
$ws.Range("J132").Value = 2265.6667
$ws.Range("K132").Value = 7200
$ws.Range("L132").Value = 20391.0003
$ws.Range("M132").Value = -4670
$ws.Range("N132").Value = -25451.0003

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H122").Value = 501574.7
$ws.Range("I122").Value = 557138.5600000001
$ws.Range("K122").Value = 1671415.68
$ws.Range("M122").Value = -1668965.68

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 826.3333
$ws.Range("I22").Value = 1280
$ws.Range("K22").Value = 1280
$ws.Range("M22").Value = -985
$ws.Range("H27").Value = 826.3333
$ws.Range("I27").Value = 1280
$ws.Range("K27").Value = 1280
$ws.Range("M27").Value = -1173
$ws.Range("H46").Value = 5320.1904
$ws.Range("I46").Value = 4818.4287
$ws.Range("J46").Value = 5571.0713
$ws.Range("K46").Value = 4818.4287
$ws.Range("L46").Value = 5571.0713
$ws.Range("M46").Value = -4630.4287
$ws.Range("N46").Value = -5947.0713
$ws.Range("H53").Value = 14450
$ws.Range("I53").Value = 14450
$ws.Range("K53").Value = 14450
$ws.Range("M53").Value = -13932
$ws.Range("H68").Value = 9209.923000000001
$ws.Range("I68").Value = 8416.333000000001
$ws.Range("K68").Value = 8416.333000000001
$ws.Range("M68").Value = -7667.333000000001
$ws.Range("H71").Value = 9209.923000000001
$ws.Range("I71").Value = 8416.333000000001
$ws.Range("K71").Value = 42081.665
$ws.Range("M71").Value = -38337.665
$ws.Range("H82").Value = 3391.0715
$ws.Range("J82").Value = 4197.375
$ws.Range("L82").Value = 4197.375
$ws.Range("N82").Value = -4919.375
$ws.Range("H85").Value = 3391.0715
$ws.Range("J85").Value = 4197.375
$ws.Range("L85").Value = 4197.375
$ws.Range("N85").Value = -6693.375
$ws.Range("H93").Value = 1310.8182
$ws.Range("I93").Value = 1291.9
$ws.Range("K93").Value = 1291.9
$ws.Range("M93").Value = -43.90000000000009
$ws.Range("H99").Value = 20000
$ws.Range("I99").Value = 20000
$ws.Range("K99").Value = 20000
$ws.Range("M99").Value = -17005
$ws.Range("H100").Value = 7579.9
$ws.Range("I100").Value = 2899
$ws.Range("J100").Value = 8100
$ws.Range("K100").Value = 2899
$ws.Range("L100").Value = 8100
$ws.Range("M100").Value = -2358
$ws.Range("N100").Value = -9182

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 11997
$ws.Range("J62").Value = 11997
$ws.Range("L62").Value = 11997
$ws.Range("N62").Value = -13245
$ws.Range("H65").Value = 11997
$ws.Range("J65").Value = 11997
$ws.Range("L65").Value = 59985
$ws.Range("N65").Value = -66225
$ws.Range("H81").Value = 900.6
$ws.Range("I81").Value = 875.25
$ws.Range("K81").Value = 1750.5
$ws.Range("M81").Value = -689.5
$ws.Range("H84").Value = 900.6
$ws.Range("I84").Value = 875.25
$ws.Range("K84").Value = 8752.5
$ws.Range("M84").Value = -3448.5
$ws.Range("H100").Value = 4500.5
$ws.Range("I100").Value = 4500.5
$ws.Range("K100").Value = 9001
$ws.Range("M100").Value = -8460
$ws.Range("H122").Value = 2295
$ws.Range("I122").Value = 2243.75
$ws.Range("K122").Value = 6731.25
$ws.Range("M122").Value = -4281.25
$ws.Range("H136").Value = 2985.5
$ws.Range("J136").Value = 9616.666999999999
$ws.Range("L136").Value = 28850.001
$ws.Range("N136").Value = -33950.001

